# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 685
$ws1.Range("F5").Value = 24
$ws1.Range("F7").Value = 47
$ws1.Range("F8").Value = 3337
$ws1.Range("F9").Value = 4269
$ws1.Range("F10").Value = 120

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 63

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 685
$ws4.Range("F5").Value = 24
$ws4.Range("F7").Value = 47
$ws4.Range("F8").Value = 3337
$ws4.Range("F9").Value = 4269
$ws4.Range("F10").Value = 120
$ws4.Range("F11").Value = 63
